$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "90.532.55"
$ws.Range("E2").Value = "  +0.16%  "

# Row 3
$ws.Range("D3").Value = "3.097.21"
$ws.Range("E3").Value = "  -1.56%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.20%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +9.17%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "625.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.62%  "

# Row 7
$ws.Range("E7").Value = "  -1.86%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.365"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.24%  "

# Row 9
$ws.Range("E9").Value = "  -0.01%  "

# Row 10
$ws.Range("D10").Value = "3.094.50"
$ws.Range("E10").Value = "  -1.62%  "

# Row 11
$ws.Range("E11").Value = "  -3.83%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.197"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.31%  "

# Row 13
$ws.Range("B13").Value = "ShibaInu"
$ws.Range("C13").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000256"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.64%  "

# Row 14
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.67"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.08%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.47"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.88%  "

# Row 16
$ws.Range("D16").Value = "90.082.77"
$ws.Range("E16").Value = "  -0.35%  "

# Row 17
$ws.Range("E17").Value = "  -2.14%  "

# Row 18
$ws.Range("D18").Value = "3.081.87"
$ws.Range("E18").Value = "  -2.94%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.79"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.94%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000217"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.70%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.02"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.34%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "436.92"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.31%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.56"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.68%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.91"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.31%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.95"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.65%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.59"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.15%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "88.73"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.65%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.23"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.92%  "

# Row 29
$ws.Range("E29").Value = "  -2.10%  "

# Row 30
$ws.Range("E30").Value = "  +0.03%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.46"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.42%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.159"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.68%  "

# Row 33
$ws.Range("E33").Value = "  +11.48%  "

# Row 34
$ws.Range("B34").Value = "dogwifhat"
$ws.Range("C34").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.88"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.05%  "

# Row 35
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.154"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.85%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "25.85"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.32%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "507.96"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.92%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.15"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.65%  "

# Row 40
$ws.Range("E40").Value = "  +0.48%  "

# Row 41
$ws.Range("E41").Value = "  +1.60%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0873"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.73%  "

# Row 43
$ws.Range("B43").Value = "PolygonEcosystemToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.411"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.87%  "

# Row 44
$ws.Range("B44").Value = "WhiteBITCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.19"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.06%  "

# Row 46
$ws.Range("E46").Value = "  +53.32%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.91"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.19%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "151.26"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.96%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.689"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.14%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "44.94"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.05%  "

# Row 51
$ws.Range("B51").Value = "FLOKI"
$ws.Range("C51").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.000277"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +17.18%  "
